# eDNA workflow with 03 and 04
# Swap the ASV_ID/Species_name/Common_name/Category values between
# row 42 <-> row 43, and row 55 <-> row 56 (columns A:D only).
# ASV_sum (E) and ASV_rank (F) stay with their original row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-RowData($Worksheet, $Row1, $Row2) {
    # Columns A through D (1 through 4)
    for ($col = 1; $col -le 4; $col++) {
        $cell1 = $Worksheet.Cells.Item($Row1, $col)
        $cell2 = $Worksheet.Cells.Item($Row2, $col)

        $value1 = $cell1.Value2
        $value2 = $cell2.Value2

        $cell1.Value2 = $value2
        $cell2.Value2 = $value1
    }
}

Swap-RowData $ws 42 43
Swap-RowData $ws 55 56
